# Refresh the cryptocurrency price / 1h-volume table with the latest
# scraped values (GitHub Actions scheduled update).
#
# Price cells that look like plain numbers are written with a leading
# apostrophe so Excel stores them as text (preserving formats such as
# trailing zeros) instead of silently converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '69.852.65'
    'E2' = '  +0.94%  '
    'D3' = '3.943.28'
    'E3' = '  +2.33%  '
    'D4' = '''0.999'
    'E4' = '  -0.03%  '
    'D5' = '''609.19'
    'E5' = '  +1.05%  '
    'D6' = '''169.62'
    'E6' = '  +3.10%  '
    'D7' = '3.940.71'
    'E7' = '  +2.40%  '
    'E8' = '  +0.00%  '
    'E9' = '  +0.36%  '
    'E10' = '  +1.75%  '
    'D11' = '''6.49'
    'E11' = '  +2.32%  '
    'E12' = '  +1.34%  '
    'E13' = '  +4.96%  '
    'D14' = '''38.10'
    'E14' = '  +2.53%  '
    'D15' = '4.605.19'
    'E15' = '  +2.39%  '
    'D16' = '3.936.98'
    'E16' = '  +1.11%  '
    'D17' = '69.896.14'
    'E17' = '  +0.84%  '
    'E18' = '  -0.72%  '
    'D19' = '''17.57'
    'E19' = '  +2.06%  '
    'E20' = '  -1.83%  '
    'D21' = '''11.07'
    'E21' = '  -5.05%  '
    'D22' = '''499.55'
    'E22' = '  +2.29%  '
    'D23' = '''0.737'
    'E23' = '  +1.89%  '
    'E24' = '  +5.63%  '
    'D25' = '''85.53'
    'E25' = '  +1.12%  '
    'D26' = '''2.29'
    'E26' = '  +1.14%  '
    'E27' = '  +0.35%  '
    'D28' = '''10.27'
    'E28' = '  +2.42%  '
    'E30' = '  +0.35%  '
    'D31' = '4.097.65'
    'E31' = '  +2.22%  '
    'E32' = '  +0.89%  '
    'D33' = '''7.89'
    'E33' = '  -1.05%  '
    'D34' = '''32.29'
    'E34' = '  -0.30%  '
    'D35' = '3.915.99'
    'E35' = '  +3.17%  '
    'E36' = '  +0.25%  '
    'B37' = 'Mantle'
    'C37' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D37' = '''1.04'
    'E37' = '  +0.83%  '
    'B38' = 'Filecoin'
    'C38' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D38' = '''6.12'
    'E38' = '  +3.45%  '
    'E39' = '  +0.58%  '
    'E40' = '  +8.81%  '
    'D41' = '''1.00'
    'E41' = '  +0.13%  '
    'E42' = '  +1.61%  '
    'E43' = '  +3.40%  '
    'D44' = '''437.37'
    'E44' = '  -0.45%  '
    'D45' = '''48.30'
    'E45' = '  -0.57%  '
    'E46' = '  +2.19%  '
    'E47' = '  +0.02%  '
    'E48' = '  +22.91%  '
    'D49' = '''0.0366'
    'E49' = '  +2.34%  '
    'D50' = '''143.21'
    'E50' = '  -0.04%  '
    'B51' = 'Arweave'
    'C51' = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
    'D51' = '''39.95'
    'E51' = '  +0.99%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
